$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-02 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-03 Sunday", 2) | Out-Null
$d.Content.Find.Execute("121÷4=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "620÷9=68, 8", 2) | Out-Null
$d.Content.Find.Execute("841÷3=280, 1", $true, $false, $false, $false, $false, $true, 1, $false, "765÷5=153, 0", 2) | Out-Null
$d.Content.Find.Execute("376÷9=41, 7", $true, $false, $false, $false, $false, $true, 1, $false, "755÷3=251, 2", 2) | Out-Null
$d.Content.Find.Execute("101÷4=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "478÷7=68, 2", 2) | Out-Null
$d.Content.Find.Execute("198÷8=24, 6", $true, $false, $false, $false, $false, $true, 1, $false, "872÷4=218, 0", 2) | Out-Null
$d.Content.Find.Execute("111÷8=13, 7", $true, $false, $false, $false, $false, $true, 1, $false, "357÷7=51, 0", 2) | Out-Null
$d.Content.Find.Execute("377÷3=125, 2", $true, $false, $false, $false, $false, $true, 1, $false, "612÷6=102, 0", 2) | Out-Null
$d.Content.Find.Execute("756÷2=378, 0", $true, $false, $false, $false, $false, $true, 1, $false, "989÷6=164, 5", 2) | Out-Null
$d.Content.Find.Execute("101÷8=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "860÷4=215, 0", 2) | Out-Null
$d.Content.Find.Execute("823÷7=117, 4", $true, $false, $false, $false, $false, $true, 1, $false, "549÷9=61, 0", 2) | Out-Null
$d.Content.Find.Execute("317÷7=45, 2", $true, $false, $false, $false, $false, $true, 1, $false, "401÷3=133, 2", 2) | Out-Null
$d.Content.Find.Execute("350÷3=116, 2", $true, $false, $false, $false, $false, $true, 1, $false, "250÷4=62, 2", 2) | Out-Null
$d.Content.Find.Execute("519÷5=103, 4", $true, $false, $false, $false, $false, $true, 1, $false, "715÷2=357, 1", 2) | Out-Null
$d.Content.Find.Execute("432÷7=61, 5", $true, $false, $false, $false, $false, $true, 1, $false, "346÷5=69, 1", 2) | Out-Null
$d.Content.Find.Execute("571÷2=285, 1", $true, $false, $false, $false, $false, $true, 1, $false, "906÷7=129, 3", 2) | Out-Null
$d.Content.Find.Execute("225÷6=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "313÷2=156, 1", 2) | Out-Null
$d.Content.Find.Execute("106÷3=35, 1", $true, $false, $false, $false, $false, $true, 1, $false, "291÷5=58, 1", 2) | Out-Null
$d.Content.Find.Execute("745÷8=93, 1", $true, $false, $false, $false, $false, $true, 1, $false, "566÷5=113, 1", 2) | Out-Null
$d.Content.Find.Execute("450÷3=150, 0", $true, $false, $false, $false, $false, $true, 1, $false, "408÷9=45, 3", 2) | Out-Null
$d.Content.Find.Execute("737÷5=147, 2", $true, $false, $false, $false, $false, $true, 1, $false, "356÷7=50, 6", 2) | Out-Null
$d.Content.Find.Execute("557÷7=79, 4", $true, $false, $false, $false, $false, $true, 1, $false, "530÷2=265, 0", 2) | Out-Null
$d.Content.Find.Execute("439÷7=62, 5", $true, $false, $false, $false, $false, $true, 1, $false, "648÷3=216, 0", 2) | Out-Null
$d.Content.Find.Execute("460÷7=65, 5", $true, $false, $false, $false, $false, $true, 1, $false, "139÷5=27, 4", 2) | Out-Null
$d.Content.Find.Execute("606÷3=202, 0", $true, $false, $false, $false, $false, $true, 1, $false, "670÷4=167, 2", 2) | Out-Null
$d.Content.Find.Execute("522÷5=104, 2", $true, $false, $false, $false, $false, $true, 1, $false, "416÷2=208, 0", 2) | Out-Null
